$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row ----
$ws.Range("A1").Value = "Nombre Problema"
$ws.Range("B1").Value = "Alpha=0"
$ws.Range("C1").Value = "Punto original (x_1,x_2,..., x_n,y_1,y_2,...,y_m)"
$ws.Range("D1").Value = "Valor Objetivo Nivel Superior"
$ws.Range("E1").Value = "Punto obtenido ahora julia (x_1,x_2,..., x_n,y_1,y_2,...,y_m)"
$ws.Range("F1").Value = "Valor Objetivo Nivel Superior Obtenido por Julia"
$ws.Range("G1").Value = "Optimizador"

# ---- Row 3: MitsosBarton2006Ex312 / TRUE ----
$ws.Range("A3").Value = "MitsosBarton2006Ex312"
$ws.Range("B3").Value = $true
$ws.Range("D3").Value = 34939
$ws.Range("D3").NumberFormat = "#,##0"

# ---- Row 4: MitsosBarton2006Ex312 / FALSE ----
$ws.Range("A4").Value = "MitsosBarton2006Ex312"
$ws.Range("B4").Value = $false
$ws.Range("D4").Value = 34939
$ws.Range("D4").NumberFormat = "#,##0"

# ---- Row 5: MitsosBarton2006Ex313 / TRUE ----
$ws.Range("A5").Value = "MitsosBarton2006Ex313"
$ws.Range("B5").Value = $true
$ws.Range("D5").Value = "   - 2.15"
$ws.Range("D5").NumberFormat = "#,##0"

# ---- Row 6: MitsosBarton2006Ex313 / FALSE ----
$ws.Range("A6").Value = "MitsosBarton2006Ex313"
$ws.Range("B6").Value = $false
$ws.Range("D6").Value = "   - 2.15"

# ---- Row 7: MitsosBarton2006Ex314 / TRUE ----
$ws.Range("A7").Value = "MitsosBarton2006Ex314"
$ws.Range("B7").Value = $true
$ws.Range("C7").Value = "(2.1,3.3)"

# ---- Row 8: MitsosBarton2006Ex314 / FALSE ----
$ws.Range("A8").Value = "MitsosBarton2006Ex314"
$ws.Range("B8").Value = $false

# ---- Row 9: MitsosBarton2006Ex323 / TRUE ----
$ws.Range("A9").Value = "MitsosBarton2006Ex323"
$ws.Range("B9").Value = $true

# ---- Row 10: MitsosBarton2006Ex323 / FALSE ----
$ws.Range("A10").Value = "MitsosBarton2006Ex323"
$ws.Range("B10").Value = $false

# ---- Row 11: MorganPatrone2006a / TRUE ----
$ws.Range("A11").Value = "MorganPatrone2006a"
$ws.Range("B11").Value = $true

# ---- Row 12: MorganPatrone2006a / FALSE ----
$ws.Range("A12").Value = "MorganPatrone2006a"
$ws.Range("B12").Value = $false

# ---- Column widths ----
$ws.Columns.Item(1).ColumnWidth = 20.0833333333333
$ws.Columns.Item(2).ColumnWidth = 27.4393939393939
$ws.Columns.Item(3).ColumnWidth = 43.4393939393939
$ws.Columns.Item(4).ColumnWidth = 26.6212121212121
$ws.Columns.Item(5).ColumnWidth = 54.8030303030303
$ws.Columns.Item(6).ColumnWidth = 44.2575757575758
$ws.Columns.Item(7).ColumnWidth = 11.6212121212121

# ---- Selection ----
$ws.Range("C6").Select()
